$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 97383937
$ws.Range("B2").Value = 77177
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 353
$ws.Range("F2").Value = 'Dvärgbägarlav'
$ws.Range("G2").Value = 'Cladonia parasitica'
$ws.Range("H2").Value = '(Hoffm.) Hoffm.'
$ws.Range("Q2").Value = 371135.0265551978
$ws.Range("R2").Value = 6744809.02867246

# Row 3
$ws.Range("A3").Value = 97383912
$ws.Range("B3").Value = 81236
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 1312
$ws.Range("F3").Value = 'Gammelgransskål'
$ws.Range("G3").Value = 'Pseudographis pinicola'
$ws.Range("H3").Value = '(Nyl.) Rehm'
$ws.Range("Q3").Value = 371300.3829165885
$ws.Range("R3").Value = 6744769.855884247

# Row 4
$ws.Range("A4").Value = 97383944
$ws.Range("B4").Value = 81236
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 1312
$ws.Range("F4").Value = 'Gammelgransskål'
$ws.Range("G4").Value = 'Pseudographis pinicola'
$ws.Range("H4").Value = '(Nyl.) Rehm'
$ws.Range("Q4").Value = 371223.8177090761
$ws.Range("R4").Value = 6744653.964307282

# Row 5
$ws.Range("A5").Value = 97383905
$ws.Range("B5").Value = 77506
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = 'Garnlav'
$ws.Range("G5").Value = 'Alectoria sarmentosa'
$ws.Range("H5").Value = '(Ach.) Ach.'
$ws.Range("Q5").Value = 371104.6060876616
$ws.Range("R5").Value = 6744765.694054871

# Row 6
$ws.Range("A6").Value = 97383922
$ws.Range("B6").Value = 78527
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 229497
$ws.Range("F6").Value = 'Korallblylav'
$ws.Range("G6").Value = 'Parmeliella triptophylla'
$ws.Range("H6").Value = '(Ach.) Müll.Arg.'
$ws.Range("Q6").Value = 371221.313728622
$ws.Range("R6").Value = 6744625.24492522

# Row 7
$ws.Range("A7").Value = 97383934
$ws.Range("B7").Value = 77595
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 6450
$ws.Range("F7").Value = 'Skuggblåslav'
$ws.Range("G7").Value = 'Hypogymnia vittata'
$ws.Range("H7").Value = '(Ach.) Parrique'
$ws.Range("Q7").Value = 371002.8556096497
$ws.Range("R7").Value = 6744724.449491068

# Row 8
$ws.Range("A8").Value = 97383928
$ws.Range("B8").Value = 78569
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 6458
$ws.Range("F8").Value = 'Lunglav'
$ws.Range("G8").Value = 'Lobaria pulmonaria'
$ws.Range("H8").Value = '(L.) Hoffm.'
$ws.Range("Q8").Value = 371053.1167658683
$ws.Range("R8").Value = 6744734.839612942

# Row 9
$ws.Range("A9").Value = 97383925
$ws.Range("B9").Value = 78569
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 6458
$ws.Range("F9").Value = 'Lunglav'
$ws.Range("G9").Value = 'Lobaria pulmonaria'
$ws.Range("H9").Value = '(L.) Hoffm.'
$ws.Range("Q9").Value = 370939.4954976452
$ws.Range("R9").Value = 6744676.446748036

# Row 10
$ws.Range("A10").Value = 97383910
$ws.Range("B10").Value = 78570
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 2081
$ws.Range("F10").Value = 'Skrovellav'
$ws.Range("G10").Value = 'Lobaria scrobiculata'
$ws.Range("H10").Value = '(Scop.) DC.'
$ws.Range("Q10").Value = 370947.5950156241
$ws.Range("R10").Value = 6744602.908073967

# Row 11
$ws.Range("A11").Value = 97383930
$ws.Range("B11").Value = 73631
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 6426
$ws.Range("F11").Value = 'Kattfotslav'
$ws.Range("G11").Value = 'Felipes leucopellaeus'
$ws.Range("H11").Value = '(Ach.) Frisch & G.Thor'
$ws.Range("Q11").Value = 371174.3428833798
$ws.Range("R11").Value = 6744597.642494702

# Row 12
$ws.Range("A12").Value = 97383908
$ws.Range("B12").Value = 78596
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 6462
$ws.Range("F12").Value = 'Stuplav'
$ws.Range("G12").Value = 'Nephroma bellum'
$ws.Range("H12").Value = '(Spreng.) Tuck.'
$ws.Range("Q12").Value = 370947.5950156241
$ws.Range("R12").Value = 6744602.908073967

# Row 13
$ws.Range("A13").Value = 97383921
$ws.Range("B13").Value = 78569
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 6458
$ws.Range("F13").Value = 'Lunglav'
$ws.Range("G13").Value = 'Lobaria pulmonaria'
$ws.Range("H13").Value = '(L.) Hoffm.'
$ws.Range("Q13").Value = 371001.0552476444
$ws.Range("R13").Value = 6744715.237031102

# Row 14
$ws.Range("A14").Value = 97383933
$ws.Range("B14").Value = 78527
$ws.Range("D14").Value = 'LC'
$ws.Range("E14").Value = 229497
$ws.Range("F14").Value = 'Korallblylav'
$ws.Range("G14").Value = 'Parmeliella triptophylla'
$ws.Range("H14").Value = '(Ach.) Müll.Arg.'
$ws.Range("Q14").Value = 371074.0237595745
$ws.Range("R14").Value = 6744569.038105329

# Row 15
$ws.Range("A15").Value = 97383909
$ws.Range("B15").Value = 78569
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 6458
$ws.Range("F15").Value = 'Lunglav'
$ws.Range("G15").Value = 'Lobaria pulmonaria'
$ws.Range("H15").Value = '(L.) Hoffm.'
$ws.Range("Q15").Value = 370947.5950156241
$ws.Range("R15").Value = 6744602.908073967

# Row 16
$ws.Range("A16").Value = 97383935
$ws.Range("B16").Value = 78569
$ws.Range("D16").Value = 'NT'
$ws.Range("E16").Value = 6458
$ws.Range("F16").Value = 'Lunglav'
$ws.Range("G16").Value = 'Lobaria pulmonaria'
$ws.Range("H16").Value = '(L.) Hoffm.'
$ws.Range("Q16").Value = 371082.9162868222
$ws.Range("R16").Value = 6744747.43465819
